$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 2 and 3: "Total Net Worth (Total Equity)" / "Positive" now come before
# "Profitable in latest Financial Statement (Latest Full Year)" / "Yes"
$ws.Range("A2").Value = "Total Net Worth (Total Equity)"
$ws.Range("B2").Value = "Positive"
$ws.Range("A3").Value = "Profitable in latest Financial Statement (Latest Full Year)"
$ws.Range("B3").Value = "Yes"

# Update Current Ratio and Gearing Ratio values
$ws.Range("B4").Value = "1.50"
$ws.Range("B5").Value = "2.00"

# Update default column width for the sheet
$ws.StandardWidth = 37.425
